$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.015.13"
$ws.Range("E2").Value = "  +1.85%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.774.83"
$ws.Range("E3").Value = "  +1.91%  "
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "328.57"
$ws.Range("E5").Value = "  +1.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4503"
$ws.Range("E7").Value = "  +0.99%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3560"
$ws.Range("E8").Value = "  +1.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07443"
$ws.Range("E9").Value = "  +1.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.95"
$ws.Range("E10").Value = "  +1.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.108"
$ws.Range("E11").Value = "  +3.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.00"
$ws.Range("E13").Value = "  +2.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.051"
$ws.Range("E14").Value = "  +2.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.262"
$ws.Range("E15").Value = "  +2.83%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.773.15"
$ws.Range("E16").Value = "  +1.60%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "94.20"
$ws.Range("E17").Value = "  +3.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001063"
$ws.Range("E18").Value = "  +1.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06443"
$ws.Range("E19").Value = "  +1.23%  "
$ws.Range("E20").Value = "  -0.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.14"
$ws.Range("E21").Value = "  +2.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.786"
$ws.Range("E22").Value = "  +1.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.016.37"
$ws.Range("E23").Value = "  +1.72%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.31"
$ws.Range("E24").Value = "  +1.90%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.132"
$ws.Range("E25").Value = "  +1.55%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.75"
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.40"
$ws.Range("E27").Value = "  +1.86%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.977.00"
$ws.Range("E28").Value = "  +1.70%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.166"
$ws.Range("E29").Value = "  +7.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.65"
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.111"
$ws.Range("E31").Value = "  +6.62%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.704"
$ws.Range("E32").Value = "  +6.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09217"
$ws.Range("E33").Value = "  +2.03%  "
$ws.Range("E34").Value = "  +1.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.88"
$ws.Range("E35").Value = "  +2.52%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06213"
$ws.Range("E36").Value = "  +3.76%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02293"
$ws.Range("E37").Value = "  +1.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2114"
$ws.Range("E38").Value = "  +2.80%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.988"
$ws.Range("E39").Value = "  +2.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6317"
$ws.Range("E40").Value = "  +1.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.188"
$ws.Range("E41").Value = "  +0.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.399"
$ws.Range("E42").Value = "  +1.57%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.914"
$ws.Range("E43").Value = "  +2.85%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.37"
$ws.Range("E44").Value = "  +2.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.756"
$ws.Range("E45").Value = "  +1.44%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5897"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.72"
$ws.Range("E47").Value = "  +0.76%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.962"
$ws.Range("E48").Value = "  +2.26%  "
$ws.Range("E49").Value = "  +2.98%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06892"
$ws.Range("E50").Value = "  +0.82%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.92"
$ws.Range("E51").Value = "  +2.57%  "
